$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Numeric value updates
$ws.Range("E2").Value = 69.42
$ws.Range("G2").Value = 635.96
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2458.59

$ws.Range("L3").Value = 56.32

$ws.Range("K4").Value = 285.12

$ws.Range("L8").Value = 1565.15

$ws.Range("L10").Value = 32.36

$ws.Range("D13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2511.66

$ws.Range("E17").Value = 69.45
$ws.Range("F17").Value = 52.25
$ws.Range("L17").Value = 1388.03

$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0

$ws.Range("L21").Value = 62.44

$ws.Range("E24").Value = 222.19

$ws.Range("L25").Value = 12043.47

# Summary row ("x de 28") updates
$ws.Range("D30").Value = "0 de 28"
$ws.Range("E30").Value = "3 de 28"
$ws.Range("F30").Value = "1 de 28"
$ws.Range("G30").Value = "1 de 28"
$ws.Range("H30").Value = "0 de 28"
$ws.Range("I30").Value = "0 de 28"
$ws.Range("J30").Value = "0 de 28"
$ws.Range("K30").Value = "1 de 28"
$ws.Range("L30").Value = "8 de 28"
